# 案件情報.xlsx — ランサーズ sheet refresh (2026-02-11 13:04 JST)
#
# New scrape pulled two extra listings ("Python SEO audit script" and the
# Salesforce PM/PL role). They land at the top of the data block (rows 7-8
# and row 9 respectively), pushing the previously-seen rows down by two.
# Every row's "取得日時" (fetched-at) timestamp is refreshed to the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Hyperlinks don't get re-targeted by a row insert in this host, so drop the
# functional links up front; they're rebuilt from the final cell text below.
$ws.Hyperlinks.Delete()

# Insert two blank rows at row 7: old rows 7-9 shift down to rows 9-11.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# Refresh the fetch timestamp on every data row (2-11).
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-02-11 13:04:39"
}

# Row 7 (new): Python SEO audit script job.
$ws.Cells.Item(7, 2).Value = "【Claude Code活用前提】PythonによるSEO監査スクリプト開発(軽量・拡張型)"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5489911"
$ws.Cells.Item(7, 7).Value = 245
$ws.Cells.Item(7, 8).Value = "🔥Python ◆開発"

# Row 8: BUYMA tool listing (was row 7 before the refresh).
$ws.Cells.Item(8, 2).Value = "【BUYMA】商品リスト取得ツールと自動出品ツール開発のご依頼"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5489608"
$ws.Cells.Item(8, 7).Value = 123
$ws.Cells.Item(8, 8).Value = "◆ツール,開発"

# Row 9 (new): Salesforce-based system build, PM/PL recruitment.
$ws.Cells.Item(9, 2).Value = "salesforceを基盤としたシステム構築案件のPM・PL募集"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5489898"
$ws.Cells.Item(9, 7).Value = 40
$ws.Cells.Item(9, 8).ClearContents()

# Row 10: domain-connection specialist wanted (was row 8; no H value).
$ws.Cells.Item(10, 2).Value = "【急募】ドメイン接続業務の専門家を探しています!"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5489674"
$ws.Cells.Item(10, 7).Value = 18
$ws.Cells.Item(10, 8).ClearContents()

# Row 11: Google Workspace setup specialist wanted (was row 9; no H value).
$ws.Cells.Item(11, 2).Value = "【急募】google work space の設定を専門家に依頼したい"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5489636"
$ws.Cells.Item(11, 7).Value = 10
$ws.Cells.Item(11, 8).ClearContents()

# Rebuild the F2:F11 hyperlinks from the (now-final) URL text in each cell.
for ($r = 2; $r -le 11; $r++) {
    $addr = $ws.Cells.Item($r, 6).Value()
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $addr) | Out-Null
}

# Column D/H got a bit wider to fit the new rows' longer text.
# (ColumnWidth round-trips through this host with a +5/6 pixel-rounding
# offset, so back it out to land on the exact stored widths 32 / 13.)
$ws.Columns.Item(4).ColumnWidth = 32 - 5/6
$ws.Columns.Item(8).ColumnWidth = 13 - 5/6
